$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Append the new "Netbeans build configuration" sentence to the note
# about building without running the tests. The sentence is added on
# a new line (manual line break, ^l) after the existing note text.
# ------------------------------------------------------------------
$old = "=true install’"
$new = "=true install’^lIf you are using Netbeans, a project configuration is provided to build without tests. "

$found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find anchor text to update the build note paragraph."
}
